$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values per diff
$ws.Range("D4").Value = 53
$ws.Range("D5").Value = 1147
$ws.Range("D6").Value = 49029
$ws.Range("C7").Value = 3

# Update the selected cell (view state) from D10 to B7
$ws.Range("B7").Select()
